$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated capital structure database
# Apply the same updated metrics to rows 2 and 3 (industry-average row and company row)

foreach ($r in 2,3) {
    $ws.Cells.Item($r, 4).Value = 0.0445                  # D: historical_growth_revenue_last_5_years
    $ws.Cells.Item($r, 5).ClearContents()                  # E: historical_growth_net_income_last_5_years (removed)

    $ws.Cells.Item($r, 7).Value = 0.148936170212766        # G: ebitdard_margin
    $ws.Cells.Item($r, 8).Value = 0.148936170212766        # H: ebitda_margin
    $ws.Cells.Item($r, 9).Value = 0.01501877346683354      # I: operating_margin
    $ws.Cells.Item($r, 10).Value = 0.01501877346683354     # J: after_tax_operating_margin
    $ws.Cells.Item($r, 11).Value = -6.06                   # K: trailing_net_income
    $ws.Cells.Item($r, 12).Value = -0.07584480600750938    # L: net_margin

    $ws.Cells.Item($r, 15).Value = 0                       # O: cash_returned_net_income
    $ws.Cells.Item($r, 18).Value = 0                       # R: payout_ratio

    $ws.Cells.Item($r, 21).Value = 16.4                    # U: cash
    $ws.Cells.Item($r, 22).Value = 0.95906432748538        # V: cash_market_cap
    $ws.Cells.Item($r, 23).Value = -0.315625                # W: roe
    $ws.Cells.Item($r, 24).Value = 0.2165775850242506       # X: cost_equity
    $ws.Cells.Item($r, 25).Value = -0.5322025850242506      # Y: roe_cost_equity
    $ws.Cells.Item($r, 26).Value = 2.15945945945946         # Z: sales_invested_capital
    $ws.Cells.Item($r, 27).Value = 0.03243243243243243      # AA: roic
    $ws.Cells.Item($r, 28).Value = 0.1118818593952615       # AB: cost_capital
    $ws.Cells.Item($r, 29).Value = -0.07944942696282904     # AC: roic_cost_capital
    $ws.Cells.Item($r, 30).Value = 27.8                     # AD: debt_total
    $ws.Cells.Item($r, 31).Value = 0                        # AE: debt_leases
    $ws.Cells.Item($r, 32).Value = 27.8                     # AF: debt_total_inc_leases
    $ws.Cells.Item($r, 33).Value = 11.4                     # AG: net_debt
    $ws.Cells.Item($r, 34).Value = 0.6191536748329621       # AH: debt_market_capital
    $ws.Cells.Item($r, 35).Value = 0.6435185185185185       # AI: debt_book_capital
    $ws.Cells.Item($r, 36).Value = 0.4                      # AJ: net_debt_market_capital
    $ws.Cells.Item($r, 37).Value = 0.4253731343283582       # AK: net_debt_book_capital
    $ws.Cells.Item($r, 38).Value = 5.13                     # AL: interest_expenses
    $ws.Cells.Item($r, 39).Value = 5.13                     # AM: net_interest_expenses
    $ws.Cells.Item($r, 40).Value = 12.69406392694064        # AN: debt_ebitda
    $ws.Cells.Item($r, 41).Value = 0.2339181286549707       # AO: ebit_interest_expenses
    $ws.Cells.Item($r, 42).Value = 5.205479452054796        # AP: net_debt_ebitda
    $ws.Cells.Item($r, 43).Value = 0.2339181286549707       # AQ: ebit_net_interest_expenses
}
